$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: season_ending_year_y 1989 -> 1988, age_y 35 -> 36
$ws.Range("Q2").Value = 1988
$ws.Range("S2").Value = 36

# Update row 3: season_ending_year_y 1998 -> 1997, age_y 25 -> 26
$ws.Range("Q3").Value = 1997
$ws.Range("S3").Value = 26
